$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 154.11111
$ws.Range("I33").Value = 154.11111
$ws.Range("K33").Value = 154.11111
$ws.Range("M33").Value = 74.88889
$ws.Range("H86").Value = 4999
$ws.Range("I86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("M86").Value = $null
$ws.Range("H89").Value = 4999
$ws.Range("I89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("M89").Value = $null
$ws.Range("H96").Value = 600.3333
$ws.Range("I96").Value = 476.57144
$ws.Range("K96").Value = 1429.71432
$ws.Range("M96").Value = -56.71432000000004
$ws.Range("H103").Value = 4398.2354
$ws.Range("I103").Value = 3212.4285
$ws.Range("J103").Value = 5228.3
$ws.Range("K103").Value = 9637.2855
$ws.Range("L103").Value = 15684.9
$ws.Range("M103").Value = -9051.2855
$ws.Range("N103").Value = -16856.9
$ws.Range("H112").Value = 2553.5557
$ws.Range("J112").Value = 2716.5
$ws.Range("L112").Value = 8149.5
$ws.Range("N112").Value = -10365.5
$ws.Range("H116").Value = 5198.625
$ws.Range("I116").Value = 4700
$ws.Range("K116").Value = 4700
$ws.Range("M116").Value = -1258
$ws.Range("H132").Value = 17689.934
$ws.Range("I132").Value = 18362.916
$ws.Range("J132").Value = 14998
$ws.Range("K132").Value = 55088.74800000001
$ws.Range("L132").Value = 44994
$ws.Range("M132").Value = -52558.74800000001
$ws.Range("N132").Value = -50054

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H47").Value = 70041
$ws.Range("J47").Value = 70041
$ws.Range("L47").Value = 70041
$ws.Range("N47").Value = -71491
$ws.Range("H74").Value = 2399.6667
$ws.Range("I74").Value = 2399.6667
$ws.Range("K74").Value = 2399.6667
$ws.Range("M74").Value = -1525.6667
$ws.Range("H77").Value = 2399.6667
$ws.Range("I77").Value = 2399.6667
$ws.Range("K77").Value = 11998.3335
$ws.Range("M77").Value = -7630.333500000001
$ws.Range("H88").Value = 2247.5
$ws.Range("I88").Value = 2837.5
$ws.Range("J88").Value = 1854.1666
$ws.Range("K88").Value = 2837.5
$ws.Range("L88").Value = 1854.1666
$ws.Range("M88").Value = -2431.5
$ws.Range("N88").Value = -2666.1666
$ws.Range("H91").Value = 2247.5
$ws.Range("I91").Value = 2837.5
$ws.Range("J91").Value = 1854.1666
$ws.Range("K91").Value = 2837.5
$ws.Range("L91").Value = 1854.1666
$ws.Range("M91").Value = -1433.5
$ws.Range("N91").Value = -4662.1666
$ws.Range("H102").Value = 6599.857
$ws.Range("I102").Value = 2066.3333
$ws.Range("K102").Value = 2066.3333
$ws.Range("M102").Value = -444.3332999999998

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H95").Value = 11601.833
$ws.Range("J95").Value = 11601.833
$ws.Range("L95").Value = 11601.833
$ws.Range("N95").Value = -17093.833
$ws.Range("H105").Value = 1671.5555
$ws.Range("I105").Value = 1435
$ws.Range("J105").Value = 2499.5
$ws.Range("K105").Value = 1435
$ws.Range("L105").Value = 2499.5
$ws.Range("M105").Value = 312
$ws.Range("N105").Value = -5993.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 90.92856999999999
$ws.Range("I7").Value = 93.90000000000001
$ws.Range("J7").Value = 83.5
$ws.Range("K7").Value = 93.90000000000001
$ws.Range("L7").Value = 83.5
$ws.Range("M7").Value = 19.09999999999999
$ws.Range("N7").Value = -309.5
$ws.Range("H16").Value = 1622.9
$ws.Range("I16").Value = 1514.3334
$ws.Range("K16").Value = 1514.3334
$ws.Range("M16").Value = -1227.3334
$ws.Range("H105").Value = 2426.5715
$ws.Range("I105").Value = 1708.4
$ws.Range("J105").Value = 4222
$ws.Range("K105").Value = 1708.4
$ws.Range("L105").Value = 4222
$ws.Range("M105").Value = 38.59999999999991
$ws.Range("N105").Value = -7716
$ws.Range("H112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("L112").Value = $null
$ws.Range("N112").Value = $null
$ws.Range("H113").Value = 1622.9
$ws.Range("I113").Value = 1514.3334
$ws.Range("K113").Value = 1514.3334
$ws.Range("M113").Value = 655.6666

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H103").Value = 1225.5714
$ws.Range("I103").Value = 1012.5
$ws.Range("J103").Value = 1310.8
$ws.Range("K103").Value = 3037.5
$ws.Range("L103").Value = 3932.4
$ws.Range("M103").Value = -2158.5
$ws.Range("N103").Value = -5690.4
$ws.Range("H114").Value = 5757.75
$ws.Range("I114").Value = 5000
$ws.Range("J114").Value = 6010.3335
$ws.Range("K114").Value = 15000
$ws.Range("L114").Value = 18031.0005
$ws.Range("M114").Value = -11746
$ws.Range("N114").Value = -24539.0005
$ws.Range("H131").Value = 966.6667
$ws.Range("I131").Value = 1000
$ws.Range("J131").Value = 960
$ws.Range("K131").Value = 3000
$ws.Range("L131").Value = 2880
$ws.Range("M131").Value = 2190
$ws.Range("N131").Value = -12960
$ws.Range("H133").Value = 4500
$ws.Range("I133").Value = 3000
$ws.Range("J133").Value = 6000
$ws.Range("K133").Value = 9000
$ws.Range("L133").Value = 18000
$ws.Range("M133").Value = -7690
$ws.Range("N133").Value = -28120
$ws.Range("H137").Value = 2300
$ws.Range("J137").Value = 2600
$ws.Range("L137").Value = 7800
$ws.Range("M137").Value = -900
$ws.Range("N137").Value = -18000

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3282.5
$ws.Range("I80").Value = 3299.6667
$ws.Range("J80").Value = 3265.3333
$ws.Range("K80").Value = 3299.6667
$ws.Range("L80").Value = 3265.3333
$ws.Range("M80").Value = -2301.6667
$ws.Range("N80").Value = -5261.3333
$ws.Range("H83").Value = 3282.5
$ws.Range("I83").Value = 3299.6667
$ws.Range("J83").Value = 3265.3333
$ws.Range("K83").Value = 16498.3335
$ws.Range("L83").Value = 16326.6665
$ws.Range("M83").Value = -11506.3335
$ws.Range("N83").Value = -26310.6665
$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = $null
$ws.Range("N101").Value = $null

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2713.5
$ws.Range("I22").Value = 675
$ws.Range("J22").Value = 3393
$ws.Range("K22").Value = 675
$ws.Range("L22").Value = 2714.75
$ws.Range("M22").Value = -380
$ws.Range("N22").Value = -3983
$ws.Range("H27").Value = 2713.5
$ws.Range("I27").Value = 675
$ws.Range("J27").Value = 3393
$ws.Range("K27").Value = 675
$ws.Range("L27").Value = 2714.75
$ws.Range("M27").Value = -568
$ws.Range("N27").Value = -3607
$ws.Range("H35").Value = 1341.5
$ws.Range("I35").Value = 1310.3334
$ws.Range("J35").Value = 1435
$ws.Range("K35").Value = 1310.3334
$ws.Range("L35").Value = 1435
$ws.Range("M35").Value = -974.3334
$ws.Range("N35").Value = -2107
$ws.Range("H40").Value = 5633.476
$ws.Range("I40").Value = 5567.4614
$ws.Range("J40").Value = 5740.75
$ws.Range("K40").Value = 5567.4614
$ws.Range("L40").Value = 5740.75
$ws.Range("M40").Value = -5431.4614
$ws.Range("N40").Value = -6012.75
$ws.Range("H47").Value = 30500
$ws.Range("I47").Value = 10000
$ws.Range("K47").Value = 10000
$ws.Range("M47").Value = -9510
$ws.Range("H52").Value = 30500
$ws.Range("I52").Value = 10000
$ws.Range("K52").Value = 10000
$ws.Range("M52").Value = -9767
$ws.Range("H69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = $null
$ws.Range("N69").Value = $null
$ws.Range("H72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = $null
$ws.Range("N72").Value = $null
$ws.Range("H82").Value = 6677.6665
$ws.Range("I82").Value = 5383.3335
$ws.Range("J82").Value = 7324.8335
$ws.Range("K82").Value = 5383.3335
$ws.Range("L82").Value = 7324.8335
$ws.Range("M82").Value = -5022.3335
$ws.Range("N82").Value = -8046.8335
$ws.Range("H85").Value = 6677.6665
$ws.Range("I85").Value = 5383.3335
$ws.Range("J85").Value = 7324.8335
$ws.Range("K85").Value = 5383.3335
$ws.Range("L85").Value = 7324.8335
$ws.Range("M85").Value = -4135.3335
$ws.Range("N85").Value = -9820.833500000001
$ws.Range("H100").Value = 5151
$ws.Range("I100").Value = 1717.1
$ws.Range("J100").Value = 8272.727999999999
$ws.Range("K100").Value = 1717.1
$ws.Range("L100").Value = 8272.727999999999
$ws.Range("M100").Value = -1176.1
$ws.Range("N100").Value = -9354.727999999999
$ws.Range("H136").Value = 4897.75
$ws.Range("I136").Value = 4897.8335
$ws.Range("J136").Value = 4897.5
$ws.Range("K136").Value = 14693.5005
$ws.Range("L136").Value = 14692.5
$ws.Range("M136").Value = -12143.5005
$ws.Range("N136").Value = -19792.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 942.1
$ws.Range("I96").Value = 912
$ws.Range("J96").Value = 962.1667
$ws.Range("K96").Value = 912
$ws.Range("L96").Value = 962.1667
$ws.Range("M96").Value = 461
$ws.Range("N96").Value = -3708.1667
$ws.Range("H136").Value = 4071.5483
$ws.Range("J136").Value = 5928.067
$ws.Range("L136").Value = 17784.201
$ws.Range("N136").Value = -22884.201
